# Interview prep.docx edit script
# Implements (per commit message "3 More Problems done" + diff):
#   1. "Problems not caught in previous 2 lists: 20" -> "...: 23", plus a new
#      paragraph "Problems Remaining: 18" right after it.
#   2. A new blank paragraph (same tab stops) inserted just above "Week 2".
#   3. "longest-palindrome" link line gets a trailing "Done" run.
#   4. "diameter-of-binary-tree" link text re-split into more runs, plus a
#      trailing " Done".
#   5. "middle-of-the-linked-list" link text re-split into more runs, plus a
#      trailing "Done" run (after the pre-existing " " run).
#   6. "01-matrix" link text re-split into more runs (no visible change).
#   7. A <w:lastRenderedPageBreak/> added to the run that starts
#      "what to study based on time left".

$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$rNs = "http://schemas.openxmlformats.org/officeDocument/2006/relationships"

# Replaces the full contents of the (unique) paragraph containing $findText
# with $innerXml, while leaving the paragraph mark itself (and therefore its
# paraId/rsid/pPr identity) untouched. $innerXml is the markup that goes
# *inside* the <w:p> (runs, hyperlinks, etc.).
function Set-ParagraphInnerXml($doc, $findText, $innerXml) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("Anchor text not found: " + $findText)
    }
    $para = $rng.Paragraphs(1)
    $start = $para.Range.Start
    $end = $para.Range.End
    $content = $doc.Range($start, $end - 1)
    $xml = '<w:p xmlns:w="' + $wNs + '" xmlns:r="' + $rNs + '">' + $innerXml + '</w:p>'
    $content.InsertXML($xml)
}

# Inserts a brand-new paragraph, described by $innerXml (+ optional $pPrXml),
# directly before the (unique) paragraph containing $findText. The found
# paragraph itself is left completely alone.
function Insert-ParagraphBefore($doc, $findText, $pPrXml, $innerXml) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("Anchor text not found: " + $findText)
    }
    $para = $rng.Paragraphs(1)
    $insPoint = $para.Range.Start
    $ins = $doc.Range($insPoint, $insPoint)
    $xml = '<w:p xmlns:w="' + $wNs + '">' + $pPrXml + $innerXml + '</w:p>'
    $ins.InsertXML($xml)
}

# Inserts a brand-new paragraph, described by $innerXml, directly after the
# (unique) paragraph containing $findText, which is otherwise left alone.
function Insert-ParagraphAfter($doc, $findText, $innerXml) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("Anchor text not found: " + $findText)
    }
    $para = $rng.Paragraphs(1)
    $insPoint = $para.Range.End - 1
    $ins = $doc.Range($insPoint, $insPoint)
    $xml = '<w:p xmlns:w="' + $wNs + '">' + $innerXml + '</w:p>'
    $ins.InsertXML($xml)
}

# --- 1. "Problems not caught ..." text fix + new "Problems Remaining" line -----
Set-ParagraphInnerXml $d "Problems not caught in previous 2 lists: 20" (
    '<w:r><w:t xml:space="preserve">Problems not caught in previous 2 lists: </w:t></w:r>' +
    '<w:r><w:t>23</w:t></w:r>'
)
Insert-ParagraphAfter $d "Problems not caught in previous 2 lists: 23" (
    '<w:r><w:t>Problems Remaining: 1</w:t></w:r>' +
    '<w:r><w:t>8</w:t></w:r>'
)

# --- 2. Blank paragraph above "Week 2" (same tab stops) ------------------------
$tabsPPr = '<w:pPr><w:tabs><w:tab w:val="left" w:pos="3204"/></w:tabs></w:pPr>'
Insert-ParagraphBefore $d "Week 2" $tabsPPr ""

# --- 3. "longest-palindrome" line: append "Done" -------------------------------
Set-ParagraphInnerXml $d "longest-palindrome" (
    '<w:hyperlink r:id="rId19" w:history="1">' +
      '<w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://leetcode.com/problems/longest-p</w:t></w:r>' +
      '<w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>a</w:t></w:r>' +
      '<w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>lindrome/</w:t></w:r>' +
    '</w:hyperlink>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>Done</w:t></w:r>'
)

# --- 4. "diameter-of-binary-tree" line: re-split hyperlink runs + " Done" ------
Set-ParagraphInnerXml $d "diameter-of-binary-tree" (
    '<w:hyperlink r:id="rId21" w:history="1">' +
      '<w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://leetcode.com/problems/diameter-o</w:t></w:r>' +
      '<w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>f</w:t></w:r>' +
      '<w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>-bi</w:t></w:r>' +
      '<w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>n</w:t></w:r>' +
      '<w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>ary-tree/</w:t></w:r>' +
    '</w:hyperlink>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>Done</w:t></w:r>'
)

# --- 5. "middle-of-the-linked-list" line: re-split hyperlink runs + "Done" -----
Set-ParagraphInnerXml $d "middle-of-the-linked-list" (
    '<w:hyperlink r:id="rId22" w:history="1">' +
      '<w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://leetcode.com/problems/middle-of-th</w:t></w:r>' +
      '<w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>e</w:t></w:r>' +
      '<w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>-linked-list/</w:t></w:r>' +
    '</w:hyperlink>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>Done</w:t></w:r>'
)

# --- 6. "01-matrix" line: re-split hyperlink runs (text unchanged) ------------
Set-ParagraphInnerXml $d "01-matrix" (
    '<w:hyperlink r:id="rId25" w:history="1">' +
      '<w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://leetcode.com/prob</w:t></w:r>' +
      '<w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>l</w:t></w:r>' +
      '<w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>ems/01-matrix/</w:t></w:r>' +
    '</w:hyperlink>'
)

# --- 7. "what to study based on time left": add lastRenderedPageBreak ---------
Set-ParagraphInnerXml $d "what to study based on time left" (
    '<w:r><w:lastRenderedPageBreak/><w:t>what to study based on time left</w:t></w:r>'
)

Write-Output "done"
